# Refresh the workers-ranking table: realeffort scores (col G) were recomputed,
# which both changes the scores themselves and re-sorts a couple of
# adjacent rows within each race group (Asian rows 2-13, Hispanic rows 14-25)
# since the sheet is kept sorted descending by realeffort within each group.
# Where two rows swap order, their prolificid (D), name (E), gender (F, when it
# differs) and index (C) travel with the person, while A/B/H/I (row counters
# and race) stay put.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 11.45740717551576
$ws.Range("G3").Value = 10.08421220545006
$ws.Range("G4").Value = 8.474440037169567
$ws.Range("G5").Value = 8.469079889133782
$ws.Range("G6").Value = 7.167398355129854

# rows 7 & 8 swap (Tu <-> Roshni)
$ws.Range("C7").Value = 16
$ws.Range("D7").Value = "60863a15760523386e761cfb"
$ws.Range("E7").Value = "Roshni"
$ws.Range("F7").Value = "female"
$ws.Range("G7").Value = 6.194924391488313

$ws.Range("C8").Value = 13
$ws.Range("D8").Value = "5697d4ae7183b8000d0fc201"
$ws.Range("E8").Value = "Tu"
$ws.Range("F8").Value = "male"
$ws.Range("G8").Value = 6.053014121671316

# rows 9 & 10 swap (Annes <-> Ankai)
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = "5c27de12a2b00a00018b2c16"
$ws.Range("E9").Value = "Ankai"
$ws.Range("F9").Value = "male"
$ws.Range("G9").Value = 5.405704693608066

$ws.Range("C10").Value = 3
$ws.Range("D10").Value = "60bd88b8fc436774352f53b9"
$ws.Range("E10").Value = "Annes"
$ws.Range("F10").Value = "female"
$ws.Range("G10").Value = 5.190617474046819

$ws.Range("G11").Value = 4.0828279198773
$ws.Range("G12").Value = 2.154970194355574
$ws.Range("G13").Value = 0.2602112850569155

# rows 14 & 15 swap (Melissa <-> Katherine); gender is female on both, unchanged
$ws.Range("C14").Value = 7
$ws.Range("D14").Value = "6024c18b094ac71dd93f4f5a"
$ws.Range("E14").Value = "Katherine"
$ws.Range("G14").Value = 8.316648944792245

$ws.Range("C15").Value = 2
$ws.Range("D15").Value = "60778ed0fde3e9c3a96f1d11"
$ws.Range("E15").Value = "Melissa"
$ws.Range("G15").Value = 8.023344841524992

$ws.Range("G16").Value = 7.429443214079729
$ws.Range("G17").Value = 7.094459853851288

# rows 18 & 19 swap (Carlos <-> Yonifredy); gender is male on both, unchanged
$ws.Range("C18").Value = 0
$ws.Range("D18").Value = "5eeaa065c7acf61c4322f6d9"
$ws.Range("E18").Value = "Yonifredy"
$ws.Range("G18").Value = 6.304187637973969

$ws.Range("C19").Value = 11
$ws.Range("D19").Value = "5f5ea8227fa75676f56f9276"
$ws.Range("E19").Value = "Carlos"
$ws.Range("G19").Value = 6.19077550199683

$ws.Range("G20").Value = 5.065205973220809

# rows 21 & 22 swap (Maria <-> Mary); gender is female on both, unchanged
$ws.Range("C21").Value = 1
$ws.Range("D21").Value = "5e0adc8f4cac6834756db412"
$ws.Range("E21").Value = "Mary"
$ws.Range("G21").Value = 3.358111939047832

$ws.Range("C22").Value = 4
$ws.Range("D22").Value = "5e706891c396cc64388ef760"
$ws.Range("E22").Value = "Maria"
$ws.Range("G22").Value = 3.344821734808749

$ws.Range("G23").Value = 2.475630392065158
$ws.Range("G24").Value = 1.281029176420817
$ws.Range("G25").Value = 0.379747773547242
